$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("P2").Value = 98
$ws.Range("Q2").Value = 98
$ws.Range("R2").Value = 98
$ws.Range("P4").Value = 1073
$ws.Range("Q4").Value = 1073
$ws.Range("R4").Value = 1073
$ws.Range("P5").Value = 748
$ws.Range("Q5").Value = 748
$ws.Range("Q9").Value = 1105
$ws.Range("Q10").Value = 98
$ws.Range("R10").Value = 98
$ws.Range("P14").Value = 98
$ws.Range("Q16").Value = 98
$ws.Range("R16").Value = 98
$ws.Range("P17").Value = 748
$ws.Range("Q17").Value = 748
$ws.Range("R17").Value = 748
$ws.Range("P18").Value = 748
$ws.Range("P22").Value = 780
$ws.Range("R22").Value = 780
$ws.Range("P23").Value = 98
$ws.Range("Q23").Value = 98
$ws.Range("R23").Value = 98
$ws.Range("P25").Value = 1073
$ws.Range("Q25").Value = 1073
$ws.Range("R25").Value = 1073
$ws.Range("P26").Value = 748
$ws.Range("Q26").Value = 748
$ws.Range("R26").Value = 748
$ws.Range("R27").Value = 423
$ws.Range("P28").Value = 98
$ws.Range("Q28").Value = 98
$ws.Range("R28").Value = 98
$ws.Range("P29").Value = 1073
$ws.Range("Q29").Value = 1073
$ws.Range("R29").Value = 1073
$ws.Range("R31").Value = 1073
$ws.Range("P32").Value = "SF"
$ws.Range("Q32").Value = "SF"
$ws.Range("P33").Value = 423
$ws.Range("Q33").Value = 423
$ws.Range("R33").Value = 423
$ws.Range("P34").Value = 1073
$ws.Range("Q34").Value = 1073
$ws.Range("R34").Value = 1073
$ws.Range("R36").Value = 455
$ws.Range("P37").Value = "SF"
$ws.Range("Q37").Value = 780
$ws.Range("R37").Value = 1073
$ws.Range("Q38").Value = 748
$ws.Range("R38").Value = 748
$ws.Range("P39").Value = 1073
$ws.Range("Q39").Value = 1073
$ws.Range("R39").Value = 1073
$ws.Range("P40").Value = 748
$ws.Range("Q40").Value = 748
$ws.Range("R40").Value = 748
$ws.Range("Q41").Value = 98
$ws.Range("R41").Value = 423
$ws.Range("P45").Value = 1073
$ws.Range("Q45").Value = 1073
$ws.Range("R45").Value = 1073
$ws.Range("P46").Value = 748
$ws.Range("Q46").Value = 748
$ws.Range("R46").Value = 748
$ws.Range("P47").Value = 423
$ws.Range("Q47").Value = 423
$ws.Range("R47").Value = 423
$ws.Range("P49").Value = 1073
$ws.Range("Q49").Value = 1073
$ws.Range("R49").Value = 1073
$ws.Range("P50").Value = 423
$ws.Range("Q50").Value = 423
$ws.Range("R50").Value = 423
$ws.Range("P51").Value = 748
$ws.Range("Q51").Value = 748
$ws.Range("R51").Value = 748
$ws.Range("P52").Value = 423
$ws.Range("Q52").Value = 423
$ws.Range("R52").Value = 423
$ws.Range("P55").Value = 1073
$ws.Range("Q55").Value = 1073
$ws.Range("R55").Value = 1073
$ws.Range("P56").Value = 748
$ws.Range("Q56").Value = 748
$ws.Range("R56").Value = 748
$ws.Range("P57").Value = 1073
$ws.Range("Q57").Value = 1073
$ws.Range("P59").Value = 748
$ws.Range("Q59").Value = 748
$ws.Range("R59").Value = 748
$ws.Range("P60").Value = 455
$ws.Range("Q60").Value = 423
$ws.Range("R60").Value = 423
$ws.Range("R62").Value = 130
$ws.Range("Q63").Value = 98
$ws.Range("P64").Value = 423
$ws.Range("Q64").Value = 98
$ws.Range("R64").Value = 423
$ws.Range("P66").Value = 1073
$ws.Range("Q66").Value = 1073
$ws.Range("R66").Value = 1073
$ws.Range("P67").Value = 748
$ws.Range("Q67").Value = 748
$ws.Range("R67").Value = 748
$ws.Range("P71").Value = 748
$ws.Range("Q71").Value = 748
$ws.Range("R71").Value = 748
$ws.Range("Q72").Value = 1073
$ws.Range("P73").Value = 423
$ws.Range("Q73").Value = 423
$ws.Range("R73").Value = 423
$ws.Range("P74").Value = 423
$ws.Range("Q74").Value = 423
$ws.Range("Q75").Value = 98
$ws.Range("R75").Value = 98
$ws.Range("P77").Value = 98
$ws.Range("Q77").Value = 98
$ws.Range("R77").Value = 98
$ws.Range("P78").Value = 1073
$ws.Range("Q78").Value = 1073
$ws.Range("R78").Value = 1073
$ws.Range("P79").Value = 748
$ws.Range("R79").Value = 423
$ws.Range("Q80").Value = 1073
$ws.Range("R82").Value = 423
$ws.Range("P84").Value = 1073
$ws.Range("P86").Value = 98
$ws.Range("Q86").Value = 98
$ws.Range("R86").Value = 98
$ws.Range("P89").Value = 98
$ws.Range("Q89").Value = 98
$ws.Range("R90").Value = 423
$ws.Range("P91").Value = 780
$ws.Range("R92").Value = 98
$ws.Range("P94").Value = 423
$ws.Range("Q94").Value = 423
$ws.Range("R94").Value = 423
$ws.Range("Q95").Value = 1073
$ws.Range("R95").Value = 1073
$ws.Range("P97").Value = 423
$ws.Range("P98").Value = 98
$ws.Range("Q98").Value = 98
$ws.Range("R98").Value = 98
$ws.Range("Q99").Value = 98
$ws.Range("R99").Value = 98
